$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1318163244425179
$ws.Range("D2").Value = 0.07493699281090471
$ws.Range("E2").Value = 0.1184991365948278
$ws.Range("F2").Value = 2.14739650687585
$ws.Range("G2").Value = 2.220330774798867
$ws.Range("H2").Value = 1.617759182396071
$ws.Range("I2").Value = 1.971795632940967
$ws.Range("L2").Value = 0.1863415131911523
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.1321083045426263
$ws.Range("D3").Value = 0.07256598175449369
$ws.Range("E3").Value = 0.1166775667030038
$ws.Range("F3").Value = 2.03681669932277
$ws.Range("G3").Value = 2.086944517262197
$ws.Range("H3").Value = 1.563157948875727
$ws.Range("I3").Value = 1.882514759531929
$ws.Range("L3").Value = 0.1808578132491476
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.1323379769755633
$ws.Range("D4").Value = 0.0711647553680308
$ws.Range("E4").Value = 0.1156327084657427
$ws.Range("F4").Value = 1.970476765595805
$ws.Range("G4").Value = 2.006664097394292
$ws.Range("H4").Value = 1.530759656186461
$ws.Range("I4").Value = 1.829116384459283
$ws.Range("L4").Value = 0.177621056856033
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.1324442345074388
$ws.Range("D5").Value = 0.07060729053340964
$ws.Range("E5").Value = 0.115225296598652
$ws.Range("F5").Value = 1.943828266499054
$ws.Range("G5").Value = 1.974349319995582
$ws.Range("H5").Value = 1.517837414266324
$ws.Range("I5").Value = 1.807708437273959
$ws.Range("L5").Value = 0.1763344579141659
$ws.Range("N5").Value = 2.293303068605894
$ws.Range("C6").Value = 0.1324626431518823
$ws.Range("D6").Value = 0.07051553691058388
$ws.Range("E6").Value = 0.1151587523329951
$ws.Range("F6").Value = 1.939426421979988
$ws.Range("G6").Value = 1.969007440302136
$ws.Range("H6").Value = 1.515708527668693
$ws.Range("I6").Value = 1.804174803194385
$ws.Range("L6").Value = 0.176122766289744
$ws.Range("N6").Value = 2.2803432614038
$ws.Range("C7").Value = 0.1323393587342636
$ws.Range("D7").Value = 0.07115718258791048
$ws.Range("E7").Value = 0.1156271397298738
$ws.Range("F7").Value = 1.970115821120544
$ws.Range("G7").Value = 2.006226678036995
$ws.Range("H7").Value = 1.530584251319738
$ws.Range("I7").Value = 1.828826249198585
$ws.Range("L7").Value = 0.1776035745791589
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.1319065342164976
$ws.Range("D8").Value = 0.07410802949213036
$ws.Range("E8").Value = 0.1178557056105518
$ws.Range("F8").Value = 2.108941954639647
$ws.Range("G8").Value = 2.173998797862311
$ws.Range("H8").Value = 1.5986968533891
$ws.Range("I8").Value = 1.94071379485068
$ws.Range("L8").Value = 0.1844234845839736
$ws.Range("N8").Value = 2.766433886209882
$ws.Range("C9").Value = 0.1314581354307478
$ws.Range("D9").Value = 0.0803363272250408
$ws.Range("E9").Value = 0.1228161569043102
$ws.Range("F9").Value = 2.393812729005703
$ws.Range("G9").Value = 2.516187907862218
$ws.Range("H9").Value = 1.741358681240058
$ws.Range("I9").Value = 2.171631052034115
$ws.Range("L9").Value = 0.198847723056204
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.131373666256188
$ws.Range("D10").Value = 0.08519417982111577
$ws.Range("E10").Value = 0.1268298777829955
$ws.Range("F10").Value = 2.611229275406941
$ws.Range("G10").Value = 2.776138822976918
$ws.Range("H10").Value = 1.851939957587604
$ws.Range("I10").Value = 2.348650619898422
$ws.Range("L10").Value = 0.2101106321087798
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.1313886564268358
$ws.Range("D11").Value = 0.08746820656500631
$ws.Range("E11").Value = 0.1287381493571544
$ws.Range("F11").Value = 2.71199586146767
$ws.Range("G11").Value = 2.896365912390365
$ws.Range("H11").Value = 1.903549187442422
$ws.Range("I11").Value = 2.430859379207675
$ws.Range("L11").Value = 0.2153845554648939
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.1314020319691309
$ws.Range("D12").Value = 0.08833877651731825
$ws.Range("E12").Value = 0.1294727822322486
$ws.Range("F12").Value = 2.750428839540348
$ws.Range("G12").Value = 2.942185579605109
$ws.Range("H12").Value = 1.923283938611632
$ws.Range("I12").Value = 2.462237726719025
$ws.Range("L12").Value = 0.2174037231214214
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("C13").Value = 0.1313988085885285
$ws.Range("D13").Value = 0.08815086007554385
$ws.Range("E13").Value = 0.1293140290502137
$ws.Range("F13").Value = 2.742139272274102
$ws.Range("G13").Value = 2.93230435937727
$ws.Range("H13").Value = 1.919025130824764
$ws.Range("I13").Value = 2.455468728007446
$ws.Range("L13").Value = 0.2169678718059487
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.1313896024319234
$ws.Range("D14").Value = 0.08753963830743317
$ws.Range("E14").Value = 0.1287983463696989
$ws.Range("F14").Value = 2.71515221639828
$ws.Range("G14").Value = 2.900129617459527
$ws.Range("H14").Value = 1.905168918343975
$ws.Range("I14").Value = 2.433435896936317
$ws.Range("L14").Value = 0.2155502292710594
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("C15").Value = 0.1313849665994837
$ws.Range("D15").Value = 0.08716648397576421
$ws.Range("E15").Value = 0.1284840450033826
$ws.Range("F15").Value = 2.698657870953383
$ws.Range("G15").Value = 2.880459975010581
$ws.Range("H15").Value = 1.896706634629652
$ws.Range("I15").Value = 2.41997258242634
$ws.Range("L15").Value = 0.2146847662822466
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.1313737629808287
$ws.Range("D16").Value = 0.08504687654055942
$ws.Range("E16").Value = 0.1267068397392741
$ws.Range("F16").Value = 2.604682019297172
$ws.Range("G16").Value = 2.768322137877874
$ws.Range("H16").Value = 1.848593757447588
$ws.Range("I16").Value = 2.343312420629076
$ws.Range("L16").Value = 0.2097690306656119
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("C17").Value = 0.1313805833460862
$ws.Range("D17").Value = 0.08376315395827305
$ws.Range("E17").Value = 0.1256378024306883
$ws.Range("F17").Value = 2.547512892194845
$ws.Range("G17").Value = 2.700040740622342
$ws.Range("H17").Value = 1.81941501844409
$ws.Range("I17").Value = 2.296718729053538
$ws.Range("L17").Value = 0.2067921978942309
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.1313895328188295
$ws.Range("D18").Value = 0.08303081680871571
$ws.Range("E18").Value = 0.125030667322676
$ws.Range("F18").Value = 2.514805620256681
$ws.Range("G18").Value = 2.660952385379346
$ws.Range("H18").Value = 1.80275480527331
$ws.Range("I18").Value = 2.270077198096743
$ws.Range("L18").Value = 0.2050941298816582
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("C19").Value = 0.1313934257489962
$ws.Range("D19").Value = 0.08278388846079565
$ws.Range("E19").Value = 0.1248264274964477
$ws.Range("F19").Value = 2.503761328198323
$ws.Range("G19").Value = 2.647749312344217
$ws.Range("H19").Value = 1.797134891359406
$ws.Range("I19").Value = 2.261083781789097
$ws.Range("L19").Value = 0.2045216068815421
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("C20").Value = 0.1313793369484699
$ws.Range("D20").Value = 0.08389918317325851
$ws.Range("E20").Value = 0.1257508003500547
$ws.Range("F20").Value = 2.553580491914317
$ws.Range("G20").Value = 2.707290174009074
$ws.Range("H20").Value = 1.822508426783827
$ws.Range("I20").Value = 2.301662320176661
$ws.Range("L20").Value = 0.2071076216461449
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("C21").Value = 0.1313920974015872
$ws.Range("D21").Value = 0.087718910842284
$ws.Range("E21").Value = 0.1289494875663806
$ws.Range("F21").Value = 2.723071452625533
$ws.Range("G21").Value = 2.909572118790948
$ws.Range("H21").Value = 1.909233597299249
$ws.Range("I21").Value = 2.4399007043246
$ws.Range("L21").Value = 0.2159660234533618
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("C22").Value = 0.131445322553823
$ws.Range("D22").Value = 0.09027048411050487
$ws.Range("E22").Value = 0.1311101069873573
$ws.Range("F22").Value = 2.835449054683608
$ws.Range("G22").Value = 3.043482725880779
$ws.Range("H22").Value = 1.967031277234753
$ws.Range("I22").Value = 2.531693937023761
$ws.Range("L22").Value = 0.2218842146923095
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("C23").Value = 0.1314128019892706
$ws.Range("D23").Value = 0.08890354055613159
$ws.Range("E23").Value = 0.129950475753553
$ws.Range("F23").Value = 2.775321754751843
$ws.Range("G23").Value = 2.971853034441551
$ws.Range("H23").Value = 1.936080023239697
$ws.Range("I23").Value = 2.482567838426519
$ws.Range("L23").Value = 0.2187136442774289
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.1313798847847423
$ws.Range("D24").Value = 0.08383766665208725
$ws.Range("E24").Value = 0.125699690739598
$ws.Range("F24").Value = 2.550836832179698
$ws.Range("G24").Value = 2.704012183990926
$ws.Range("H24").Value = 1.821109538856319
$ws.Range("I24").Value = 2.299426868926901
$ws.Range("L24").Value = 0.2069649770088091
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.1315364887495889
$ws.Range("D25").Value = 0.07860283156433923
$ws.Range("E25").Value = 0.1214100718260696
$ws.Range("F25").Value = 2.315352257043571
$ws.Range("G25").Value = 2.422153262617257
$ws.Range("H25").Value = 1.70176852843656
$ws.Range("I25").Value = 2.107894272092096
$ws.Range("L25").Value = 0.1948305746392549
$ws.Range("N25").Value = 3.331249627311138
